{"js": "// Replace the trailing empty paragraph (just before the section break) with\n// two new discussion paragraphs describing future TTE-modeling work and a\n// comparison with energy-development research, matching the author's commit.\n//\n// We build the exact WordprocessingML for the two paragraphs (including the\n// w:proofErr spell/grammar-check bookmarks the original Word session left\n// behind) and insert it via Paragraph.insertOoxml (Flat OPC wrapped), which\n// lets us control run boundaries and formatting precisely instead of\n// relying on auto-split plain-text insertion.\n\nconst newParagraphsXml =\n  '<w:p>' +\n    '<w:pPr>' +\n      '<w:rPr>' +\n        '<w:rFonts w:eastAsia=\"Times New Roman\"/>' +\n        '<w:color w:val=\"000000\"/>' +\n      '</w:rPr>' +\n    '</w:pPr>' +\n    '<w:r>' +\n      '<w:rPr>' +\n        '<w:rFonts w:eastAsia=\"Times New Roman\"/>' +\n        '<w:color w:val=\"000000\"/>' +\n      '</w:rPr>' +\n      '<w:t>Further study could build on our research using TTE modeling as a less intensive and invasive method for estimating density (Moeller et al. 2018</w:t>' +\n    '</w:r>' +\n    '<w:r>' +\n      '<w:rPr>' +\n        '<w:rFonts w:eastAsia=\"Times New Roman\"/>' +\n        '<w:color w:val=\"000000\"/>' +\n      '</w:rPr>' +\n      '<w:t xml:space=\"preserve\">, </w:t>' +\n    '</w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r>' +\n      '<w:rPr>' +\n        '<w:rFonts w:eastAsia=\"Times New Roman\"/>' +\n        '<w:color w:val=\"000000\"/>' +\n      '</w:rPr>' +\n      '<w:t>Loonam</w:t>' +\n    '</w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r>' +\n      '<w:rPr>' +\n        '<w:rFonts w:eastAsia=\"Times New Roman\"/>' +\n        '<w:color w:val=\"000000\"/>' +\n      '</w:rPr>' +\n      '<w:t xml:space=\"preserve\"> et al. 2021), while following further development of TTE study design (Moeller et al. 2023). Maximizing data derived from collared individuals including survival, reproduction and nutritional condition can help identify nutritional carrying capacity, revealing the role recreation plays in the limitation of wildlife species and other density-independent impacts recreation might have at the population level. With density estimations not being particularly sensitive and having high variability from year to year collecting this supplemental data will provide support for the mechanisms we propose may be important in this population (Bergman et al. 2015). Our work represents a snapshot in time with low density possibly caused by a variety of reasons that we were not able to distinguish. Longer term study following the trends of recreation over time can resolve this uncertainty present in our study design.</w:t>' +\n    '</w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n    '<w:r>' +\n      '<w:rPr>' +\n        '<w:rFonts w:eastAsia=\"Times New Roman\"/>' +\n        '<w:color w:val=\"000000\"/>' +\n      '</w:rPr>' +\n      '<w:t xml:space=\"preserve\">Similar research into energy development did not find habituation by mule deer with respect to decreased avoidance of well pads, also observing a concomitant decline in the population over 15 years of development (Sawyer et al. 2017). </w:t>' +\n    '</w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r>' +\n      '<w:rPr>' +\n        '<w:rFonts w:eastAsia=\"Times New Roman\"/>' +\n        '<w:color w:val=\"000000\"/>' +\n      '</w:rPr>' +\n      '<w:t>However</w:t>' +\n    '</w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r>' +\n      '<w:rPr>' +\n        '<w:rFonts w:eastAsia=\"Times New Roman\"/>' +\n        '<w:color w:val=\"000000\"/>' +\n      '</w:rPr>' +\n      '<w:t xml:space=\"preserve\"> responses may differ between these two disturbance types, as the most active phase of natural gas development includes noise and artificial light beyond heavy traffic volumes and human activity (Northrup et al. 2021).</w:t>' +\n    '</w:r>' +\n  '</w:p>';\n\nconst flatOpc =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' + newParagraphsXml + '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document ends with a lone empty paragraph right before the section\n// break; that is the `<w:p/>` the diff turns into the two new paragraphs.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertOoxml(flatOpc, \"Replace\");\nawait context.sync();\n", "ps1": "# Replace the trailing empty paragraph (just before the section break) with\n# two new discussion paragraphs describing future TTE-modeling work and a\n# comparison with energy-development research, matching the author's commit.\n#\n# Range.InsertXML lets us drop in the exact WordprocessingML for the new\n# paragraphs -- including the w:proofErr spell/grammar-check bookmarks the\n# original Word session left behind -- instead of relying on plain-text\n# insertion (which would not reproduce the run boundaries / proofErr marks).\n\n$d = $word.ActiveDocument\n\n$newParagraphsXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:rPr><w:rFonts w:eastAsia=\"Times New Roman\"/><w:color w:val=\"000000\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia=\"Times New Roman\"/><w:color w:val=\"000000\"/></w:rPr><w:t>Further study could build on our research using TTE modeling as a less intensive and invasive method for estimating density (Moeller et al. 2018</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia=\"Times New Roman\"/><w:color w:val=\"000000\"/></w:rPr><w:t xml:space=\"preserve\">, </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:rFonts w:eastAsia=\"Times New Roman\"/><w:color w:val=\"000000\"/></w:rPr><w:t>Loonam</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:rFonts w:eastAsia=\"Times New Roman\"/><w:color w:val=\"000000\"/></w:rPr><w:t xml:space=\"preserve\"> et al. 2021), while following further development of TTE study design (Moeller et al. 2023). Maximizing data derived from collared individuals including survival, reproduction and nutritional condition can help identify nutritional carrying capacity, revealing the role recreation plays in the limitation of wildlife species and other density-independent impacts recreation might have at the population level. With density estimations not being particularly sensitive and having high variability from year to year collecting this supplemental data will provide support for the mechanisms we propose may be important in this population (Bergman et al. 2015). Our work represents a snapshot in time with low density possibly caused by a variety of reasons that we were not able to distinguish. Longer term study following the trends of recreation over time can resolve this uncertainty present in our study design.</w:t></w:r></w:p><w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:rPr><w:rFonts w:eastAsia=\"Times New Roman\"/><w:color w:val=\"000000\"/></w:rPr><w:t xml:space=\"preserve\">Similar research into energy development did not find habituation by mule deer with respect to decreased avoidance of well pads, also observing a concomitant decline in the population over 15 years of development (Sawyer et al. 2017). </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:rPr><w:rFonts w:eastAsia=\"Times New Roman\"/><w:color w:val=\"000000\"/></w:rPr><w:t>However</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:rPr><w:rFonts w:eastAsia=\"Times New Roman\"/><w:color w:val=\"000000\"/></w:rPr><w:t xml:space=\"preserve\"> responses may differ between these two disturbance types, as the most active phase of natural gas development includes noise and artificial light beyond heavy traffic volumes and human activity (Northrup et al. 2021).</w:t></w:r></w:p>'\n\n# The document ends with a lone empty paragraph right before the section\n# break; that is the `<w:p/>` the diff turns into the two new paragraphs.\n$target = $d.Paragraphs.Last\n$targetRange = $target.Range\n$targetRange.InsertXML($newParagraphsXml) | Out-Null\n"}
